# Generate Report for Handback
# ------------------------------------------------------------------
# This mirrors the "handback" report-generation step: the Status /
# "Ready for handoff" text becomes "Handed back: in sync with en-US"
# everywhere it appears (Overview zh-cn/de-de columns + the Status
# column on each per-language sheet), and each per-language sheet's
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" cells get populated now that the handback has completed.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdFile      = "05049e00-37bb-4c68-ae9a-126150ae4e7f.md"
$mdFileUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/18124d46eb4ca4fdc3cc1908c12bc1c245fdb246/e2e/05049e00-37bb-4c68-ae9a-126150ae4e7f.md"

$zhXlf = "05049e00-37bb-4c68-ae9a-126150ae4e7f.62faf983f671d06f1b12dade6646264544e38dfe.zh-cn.xlf"
$deXlf = "05049e00-37bb-4c68-ae9a-126150ae4e7f.62faf983f671d06f1b12dade6646264544e38dfe.de-de.xlf"

$zhHandbackTime = "2016-08-15 15:00:34"
$deHandbackTime = "2016-08-15 15:00:42"

# ------------------------------------------------------------------
# 1. Overview sheet: zh-cn / de-de status columns (E,F for rows 2,3)
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew

# Column widths widen to fit the longer status text (autofit-style).
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# ------------------------------------------------------------------
# 2. zh-cn detail sheet
# ------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

$wsZh.Range("I2").Value = $mdFile
$wsZh.Range("J2").Value = $zhXlf
$wsZh.Range("K2").Value = $zhHandbackTime

$wsZh.Range("I3").Value = $mdFile
$wsZh.Range("J3").Value = $zhXlf
$wsZh.Range("K3").Value = $zhHandbackTime

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdFileUrl, $null, $null, $mdFile)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdFileUrl, $null, $null, $mdFile)

$wsZh.Columns.Item(3).ColumnWidth = 29.14
$wsZh.Columns.Item(9).ColumnWidth = 39.1
$wsZh.Columns.Item(10).ColumnWidth = 39.1

# ------------------------------------------------------------------
# 3. de-de detail sheet
# ------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

$wsDe.Range("I2").Value = $mdFile
$wsDe.Range("J2").Value = $deXlf
$wsDe.Range("K2").Value = $deHandbackTime

$wsDe.Range("I3").Value = $mdFile
$wsDe.Range("J3").Value = $deXlf
$wsDe.Range("K3").Value = $deHandbackTime

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdFileUrl, $null, $null, $mdFile)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdFileUrl, $null, $null, $mdFile)

$wsDe.Columns.Item(3).ColumnWidth = 29.14
$wsDe.Columns.Item(9).ColumnWidth = 39.1
$wsDe.Columns.Item(10).ColumnWidth = 39.1

Write-Host "Handback report generated."
